$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '34.408.77'
$ws.Range("E2").Value = '  -0.35%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.796.98'
$ws.Range("E3").Value = '  -0.91%  '

# Row 4
$ws.Range("E4").Value = '  -0.17%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '224.19'
$ws.Range("E5").Value = '  -2.00%  '

# Row 6
$ws.Range("E6").Value = '  -1.15%  '

# Row 7
$ws.Range("E7").Value = '  -0.16%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '38.99'
$ws.Range("E8").Value = '  +6.56%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.287'
$ws.Range("E9").Value = '  -4.71%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0666'
$ws.Range("E10").Value = '  -5.49%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0986'
$ws.Range("E11").Value = '  +1.84%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.056.95'
$ws.Range("E12").Value = '  -0.99%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.794.16'
$ws.Range("E13").Value = '  -1.27%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.82'
$ws.Range("E14").Value = '  -6.02%  '

# Row 15
$ws.Range("B15").Value = 'WrappedBTC'
$ws.Range("C15").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '34.422.24'
$ws.Range("E15").Value = '  -0.35%  '

# Row 16
$ws.Range("B16").Value = 'Polygon'
$ws.Range("C16").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.627'
$ws.Range("E16").Value = '  -4.12%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.34'
$ws.Range("E17").Value = '  -3.50%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '67.81'
$ws.Range("E18").Value = '  -3.54%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '238.80'
$ws.Range("E19").Value = '  -3.15%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0763'
$ws.Range("E20").Value = '  -4.44%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.99'
$ws.Range("E21").Value = '  -5.71%  '

# Row 22
$ws.Range("E22").Value = '  -0.13%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.06'
$ws.Range("E23").Value = '  -4.01%  '

# Row 24
$ws.Range("E24").Value = '  -4.49%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '171.32'
$ws.Range("E25").Value = '  -0.45%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '17.44'
$ws.Range("E26").Value = '  -2.64%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.59'
$ws.Range("E27").Value = '  -6.22%  '

# Row 28
$ws.Range("E28").Value = '  -2.19%  '

# Row 29
$ws.Range("E29").Value = '  -0.21%  '

# Row 30
$ws.Range("E30").Value = '  -2.48%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.72'
$ws.Range("E31").Value = '  -4.03%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0512'
$ws.Range("E32").Value = '  -3.79%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.80'
$ws.Range("E33").Value = '  -5.57%  '

# Row 34
$ws.Range("E34").Value = '  -1.62%  '

# Row 35
$ws.Range("B35").Value = 'Maker'
$ws.Range("C35").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.306.30'
$ws.Range("E35").Value = '  -6.79%  '

# Row 36
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.635'
$ws.Range("E36").Value = '  -5.41%  '

# Row 37
$ws.Range("E37").Value = '  -1.88%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0184'
$ws.Range("E38").Value = '  -3.56%  '

# Row 39
$ws.Range("B39").Value = 'WEMIXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.22'
$ws.Range("E39").Value = '  +2.13%  '

# Row 40
$ws.Range("E40").Value = '  +0.65%  '

# Row 41
$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.27'
$ws.Range("E41").Value = '  -6.32%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '81.51'
$ws.Range("E42").Value = '  -1.92%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.79'
$ws.Range("E43").Value = '  -2.14%  '

# Row 44
$ws.Range("B44").Value = 'InjectiveProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.11'
$ws.Range("E44").Value = '  +2.94%  '

# Row 45
$ws.Range("B45").Value = 'ARBITRUM'
$ws.Range("C45").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.938'
$ws.Range("E45").Value = '  -3.24%  '

# Row 46
$ws.Range("E46").Value = '  +3.73%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.959.15'
$ws.Range("E47").Value = '  -0.94%  '

# Row 48
$ws.Range("B48").Value = 'PaxDollar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.00'
$ws.Range("E48").Value = '  -0.18%  '

# Row 49
$ws.Range("B49").Value = 'FraxShare'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.68'
$ws.Range("E49").Value = '  -6.35%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '101.36'
$ws.Range("E50").Value = '  -3.02%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0608'
$ws.Range("E51").Value = '  -1.33%  '
